$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header BL1 (was "Var1a_STOCKCHANGE_TOTAL") to the new label
$ws.Range("BL1").Value = "Var1_totalC_Output"

# Populate/update the BL column (rows 2-152) with the new values
$blValues = @(5861772.4768134877,11107207.699530974,12110973.317942418,13088917.898367032,14079955.044636741,15038869.737551011,16355132.155010484,17696343.05710458,17980914.267286956,19247980.76036144,20121398.878286839,20691086.233655393,21664567.954973303,22326768.209028967,22734845.238919772,23095793.93124678,24126564.668017589,24394910.870130464,24651180.615362965,25435542.248658709,26031574.055885121,25729505.444059506,26711456.699059509,27795700.508901119,28133640.416426435,28565920.361884296,28950774.756160446,29087843.84363357,29332806.683907479,29966783.942997802,29244107.445718028,28258762.940532029,27507905.550222665,27789212.677103274,27990556.903467514,28472006.006601132,29159552.032543812,29527268.53764116,29211564.667539142,29791413.629526764,30208278.488520101,31173817.404996552,31504703.182076447,31478518.687336091,31638352.776278842,31341045.332746103,32138283.587713186,32502889.082307365,32906201.019293368,32470431.400846072,33436502.743709818,36068772.551706925,36195276.529252745,36359525.320088685,36495513.217898309,36845888.368946053,37079812.52549763,36638889.516439773,36851418.178238228,37433299.215469189,37140892.939016894,37206736.84893553,37525614.063833527,38008001.577397153,38461416.73345612,39047040.816286109,39370990.825217731,39450416.506517678,39999080.230682567,40367458.604179874,40615735.003486946,41440882.010553926,42120218.771336287,42552082.222033814,42056278.70604521,41848275.269384481,42797261.404742464,43745152.40859618,44434610.698394954,44712739.419059835,44230154.241114073,44127601.905859753,44363783.162989564,45804828.860371761,46769358.424599171,47427273.101053461,48444305.917312309,49502786.917974636,49952927.46418637,50512852.148384809,50789440.409271389,50653922.480584003,51521203.174823865,52205308.524761304,52891672.850889988,53352870.286617175,54077914.362575822,54787478.003478453,55587859.252851486,56476870.788903393,57095558.678844079,57475588.403711736,58306899.410335012,58878410.969623826,60241447.722770736,61080815.078452639,61324266.776435003,60853366.689290695,59250725.680234902,58232242.619164847,58163727.003920771,58298544.497461699,58507421.308226079,58959536.607529126,59492926.333034784,60054468.394272581,60680174.483157687,61368780.056006297,62119048.182299905,62929768.898728885,64110434.056330606,65043010.815260977,65625794.718345448,66225681.965189941,66817738.292229086,67362284.731110185,67889381.776356235,68451036.827588499,68987318.217551723,69524230.496142685,70052880.927325845,70650266.539540887,71174220.481345162,71742773.793748215,72316781.767209962,72866313.371986657,73443580.676817179,73994241.330868781,74528083.338661656,75080146.843625098,75653935.344137818,76236019.02358593,76843634.577108741,77454777.414531231,78034271.763398319,78638964.061451614,79221218.7722615,79808372.30424574,80399049.186338097,80988190.030189544,81574096.161355004)
for ($i = 0; $i -lt $blValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 64).Value = $blValues[$i]
}

# Reflect the selection/view state used while editing this column
[void]$ws.Range("BL2:BL152").Select()
